# Fruta / hortaliza, semanal
# Re-shuffle the D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) values across data rows
# 2..32 according to a fixed row permutation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: destination row -> source row (values currently sitting in the source
# row are the ones that must end up in the destination row).
$map = @{
    2  = 3
    3  = 7
    4  = 15
    5  = 8
    6  = 16
    7  = 10
    8  = 31
    9  = 20
    10 = 29
    11 = 21
    12 = 2
    13 = 23
    14 = 26
    15 = 30
    16 = 13
    17 = 18
    18 = 22
    19 = 17
    20 = 25
    21 = 9
    22 = 4
    23 = 19
    24 = 24
    25 = 14
    26 = 28
    27 = 12
    28 = 32
    29 = 27
    30 = 5
    31 = 11
    32 = 6
}

$cols = @(4, 10, 11, 12, 13, 16)   # D, J, K, L, M, P

# Snapshot the original values for every affected column/row before writing
# anything back, since the permutation is not a simple swap (it has a
# 7-cycle and a 23-cycle), so a naive in-place copy would clobber sources
# that are still needed later.
$orig = @{}
for ($r = 2; $r -le 32; $r++) {
    foreach ($c in $cols) {
        $orig["$r-$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $orig["$srcRow-$c"]
    }
}
